# Prefix each step/command name in column A (row 2 onward) with the
# worksheet's own name, e.g. "Step4 Takeaway" -> "discount1 Step4 Takeaway".
# This matches the commit message: "fix: unique command names in XLSX -
# prefix protocol name to each step".

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "price1", "price2",
    "discount1", "discount2",
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol", "dickpic", "boosters"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $lastRow = $ws.UsedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $old = $cell.Value2
        if ($old -ne $null -and $old -ne "") {
            $prefix = $name + " "
            if ($old.ToString().StartsWith($prefix) -eq $false) {
                $cell.Value = $prefix + $old
            }
        }
    }
}
